$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999997662573199
$ws.Range("E2").Value = 0.9999997662573199

# Row 3
$ws.Range("D3").Value = 0.9881130429213031
$ws.Range("E3").Value = 0.9881130429213031

# Row 4
$ws.Range("D4").Value = 0.9999611323063281
$ws.Range("E4").Value = 0.9999611323063281

# Row 5
$ws.Range("D5").Value = 0.000003233928126591932
$ws.Range("E5").Value = 0.000003233928126591932

# Row 6
$ws.Range("D6").Value = 0.0000002043481267501852
$ws.Range("E6").Value = 0.0000002043481267501852

# Row 7
$ws.Range("D7").Value = 0.9999999999971303
$ws.Range("E7").Value = 0.000000000002869704474051105

# Row 8
$ws.Range("D8").Value = 0.9999899470418858
$ws.Range("E8").Value = 0.00001005295811418172

# Row 9
$ws.Range("D9").Value = 0.9999999562022234
$ws.Range("E9").Value = 0.00000004379777662766315

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.00000003514926314992929
$ws.Range("E11").Value = 0.9999999648507368
$ws.Range("F11").Value = 4.702038764953613
$ws.Range("G11").Value = 0.6
